# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# Semantic changes applied to Hoja1:
#   1. "Periodo Mora" value 2508 -> 2509 for the rows that still had it
#      (rows 20, 21 and 22 in column E).
#   2. The "Periodo Mora" column (E16:E22) is center aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1. Update the old period value (2508) to the new one (2509)
$ws.Range("E20").Value2 = "2509"
$ws.Range("E21").Value2 = "2509"
$ws.Range("E22").Value2 = "2509"

# 2. Center the "Periodo Mora" data column
$ws.Range("E16:E22").HorizontalAlignment = -4108
